$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.901.41"
$ws.Range("E2").Value = "  +2.71%  "

$ws.Range("D3").Value = "1.873.59"
$ws.Range("E3").Value = "  +1.08%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.012"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.57%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.07"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.16%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.011"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.50%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4841"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.37%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3829"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.50%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07387"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.57%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9417"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.83%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.09"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.66%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07815"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.03%  "

$ws.Range("D13").Value = "1.880.47"
$ws.Range("E13").Value = "  +1.83%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.506"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.98%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.612"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.62%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.34"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.86%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.013"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.62%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008884"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.95%  "

$ws.Range("E19").Value = "  -0.50%  "

$ws.Range("D20").Value = "28.028.06"
$ws.Range("E20").Value = "  +3.05%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.87"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.65%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.131"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.96%  "

$ws.Range("D23").Value = "2.130.46"
$ws.Range("E23").Value = "  +2.93%  "

$ws.Range("E24").Value = "  +1.73%  "

$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.26"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.63%  "

$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.932"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.63%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.59"
$ws.Range("D27").Style = "Normal"

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.064"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.71%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "116.12"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.90%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.33%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08929"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.60%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.332"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.82%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.227"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.88%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7686"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.12%  "

$ws.Range("E35").Value = "  +2.70%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.714"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.82%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.134"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.42%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02054"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.06%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5638"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.47%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05376"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.22%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.996"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.49%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.048"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.04%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.597"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.26%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1539"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.75%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4898"
$ws.Range("D45").Style = "Normal"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.68"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.08%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "105.34"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.02%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.012"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.51%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.675"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.95%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "68.15"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.11%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06116"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.90%  "
